$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --------------------------------------------------------------------------
# Add "Sheet2" right after "Sheet1". It reproduces the SAW (normalization)
# calculation from Sheet1, but the "C2" criterion (column C) is switched
# from a benefit criterion (value / max) to a cost criterion (min / value).
# --------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Weights row
$ws2.Range("B2").Value = 0.35
$ws2.Range("C2").Value = 0.25
$ws2.Range("D2").Value = 0.25
$ws2.Range("E2").Value = 0.15
$ws2.Range("B2:E2").NumberFormat = "0.000"

# Criteria header row
$ws2.Range("B3").Value = "C1"
$ws2.Range("C3").Value = "C2"
$ws2.Range("D3").Value = "C3"
$ws2.Range("E3").Value = "C4"
$ws2.Range("B3:E3").NumberFormat = "0.000"

# Raw decision-matrix rows
$ws2.Range("A4").Value = "A1"
$ws2.Range("B4").Value = 70
$ws2.Range("C4").Value = 50
$ws2.Range("D4").Value = 80
$ws2.Range("E4").Value = 60
$ws2.Range("A4:E4").NumberFormat = "0.000"

$ws2.Range("A5").Value = "A2"
$ws2.Range("B5").Value = 50
$ws2.Range("C5").Value = 60
$ws2.Range("D5").Value = 82
$ws2.Range("E5").Value = 70
$ws2.Range("A5:E5").NumberFormat = "0.000"

$ws2.Range("A6").Value = "A3"
$ws2.Range("B6").Value = 85
$ws2.Range("C6").Value = 55
$ws2.Range("D6").Value = 80
$ws2.Range("E6").Value = 75
$ws2.Range("A6:E6").NumberFormat = "0.000"

$ws2.Range("A7").Value = "A4"
$ws2.Range("B7").Value = 82
$ws2.Range("C7").Value = 70
$ws2.Range("D7").Value = 65
$ws2.Range("E7").Value = 85
$ws2.Range("A7:E7").NumberFormat = "0.000"

$ws2.Range("A8").Value = "A5"
$ws2.Range("B8").Value = 75
$ws2.Range("C8").Value = 75
$ws2.Range("D8").Value = 85
$ws2.Range("E8").Value = 74
$ws2.Range("A8:E8").NumberFormat = "0.000"

$ws2.Range("A9").Value = "A6"
$ws2.Range("B9").Value = 62
$ws2.Range("C9").Value = 50
$ws2.Range("D9").Value = 75
$ws2.Range("E9").Value = 80
$ws2.Range("A9:E9").NumberFormat = "0.000"

# Section label
$ws2.Range("B11").Value = "Hasil Normalisasi Matriks"
$ws2.Range("B11").NumberFormat = "0.000"

# Normalized matrix + weighted SAW score.
# B, D, E stay benefit criteria (value / max of column).
# C becomes a cost criterion (min of column / value) -- this is the
# "benefit -> cost" change called out in the commit message.
$ws2.Range("A12").Value = "A1 = Indra,"
$ws2.Range("B12").Formula = "=B4/MAX(B4:B9)"
$ws2.Range("C12").Formula = "=MIN(C4:C9)/C4"
$ws2.Range("D12").Formula = "=D4/MAX(D4:D9)"
$ws2.Range("E12").Formula = "=E4/MAX(E4:E9)"
$ws2.Range("F12").Formula = "=B12*B2+C12*C2+D12*D2+E12*E2"
$ws2.Range("A12:F12").NumberFormat = "0.000"

$ws2.Range("A13").Value = "A2 = Roni,"
$ws2.Range("B13").Formula = "=B5/MAX(B4:B9)"
$ws2.Range("C13").Formula = "=MIN(C4:C9)/C5"
$ws2.Range("D13").Formula = "=D5/MAX(D4:D9)"
$ws2.Range("E13").Formula = "=E5/MAX(E4:E9)"
$ws2.Range("F13").Formula = "=B13*B2+C13*C2+D13*D2+E13*E2"
$ws2.Range("A13:F13").NumberFormat = "0.000"

$ws2.Range("A14").Value = "A3 = Putri,"
$ws2.Range("B14").Formula = "=B6/MAX(B4:B9)"
$ws2.Range("C14").Formula = "=MIN(C4:C9)/C6"
$ws2.Range("D14").Formula = "=D6/MAX(D4:D9)"
$ws2.Range("E14").Formula = "=E6/MAX(E4:E9)"
$ws2.Range("F14").Formula = "=B14*B2+C14*C2+D14*D2+E14*E2"
$ws2.Range("A14:F14").NumberFormat = "0.000"

$ws2.Range("A15").Value = "A4 = Dani,"
$ws2.Range("B15").Formula = "=B7/MAX(B4:B9)"
$ws2.Range("C15").Formula = "=MIN(C4:C9)/C7"
$ws2.Range("D15").Formula = "=D7/MAX(D4:D9)"
$ws2.Range("E15").Formula = "=E7/MAX(E4:E9)"
$ws2.Range("F15").Formula = "=B15*B2+C15*C2+D15*D2+E15*E2"
$ws2.Range("A15:F15").NumberFormat = "0.000"

$ws2.Range("A16").Value = "A5 = Ratna,"
$ws2.Range("B16").Formula = "=B8/MAX(B4:B9)"
$ws2.Range("C16").Formula = "=MIN(C4:C9)/C8"
$ws2.Range("D16").Formula = "=D8/MAX(D4:D9)"
$ws2.Range("E16").Formula = "=E8/MAX(E4:E9)"
$ws2.Range("F16").Formula = "=B16*B2+C16*C2+D16*D2+E16*E2"
$ws2.Range("A16:F16").NumberFormat = "0.000"

$ws2.Range("A17").Value = "A6 = Mira. "
$ws2.Range("B17").Formula = "=B9/MAX(B4:B9)"
$ws2.Range("C17").Formula = "=MIN(C4:C9)/C9"
$ws2.Range("D17").Formula = "=D9/MAX(D4:D9)"
$ws2.Range("E17").Formula = "=E9/MAX(E4:E9)"
$ws2.Range("F17").Formula = "=B17*B2+C17*C2+D17*D2+E17*E2"
$ws2.Range("A17:F17").NumberFormat = "0.000"

# --------------------------------------------------------------------------
# View/selection bookkeeping to mirror the diff: Sheet1 is no longer the
# selected tab and its selection moves to A2:E8; Sheet2 becomes the active
# tab, zoomed to 109%, with C17 selected.
# --------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("A2:E8").Select() | Out-Null

$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 109
$ws2.Range("C17").Select() | Out-Null
